$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 already has empty C6/D6 cells (numeric, blank). The new timesheet
# entry adds the clock-in date/time into A6/B6, styled like the other
# data rows (e.g. row 5).
#
# A6 holds a date-formatted string ("2026-01-25"). Prefixing with a single
# quote forces it to be stored as literal text instead of being
# auto-converted into a date serial value.
$ws.Range("A6").Value = "'2026-01-25"
$ws.Range("B6").Value = "23:19:54"

# Match the formatting (font/border/etc.) used by the rest of the data
# rows, and drop the quote-prefix text marker picked up above so the
# cells end up styled identically to A5/B5.
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B5").Copy() | Out-Null
$ws.Range("B6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0
